# graphs/charts and diagrams slides
#
# Add links to the slide decks that were published for the
# "Methods & Procedures" (week 7) and "Charts, Tables, and Diagrams"
# (week 8) classes, and swap in the new week-8 Wednesday agenda
# ("Peer review time" instead of the IKEA-style assembly assessment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Week-plan" is the active/tab-selected sheet

# Week 7 (row 8), Wednesday_Class: link the "how to" guides write-up to its slides
$ws.Range("H8").Value = "- [Writing good “how to” guides for others](../slides/07-user-guides.qmd)"

# Week 8 (row 9), Monday_Class: link "Tables and charts" to its slides
$ws.Range("G9").Value = "- [Tables and charts](../slides/08-tables-charts.qmd)"

# Week 8 (row 9), Wednesday_Class: link "Diagrams" to its slides and
# replace the IKEA assembly-diagram assessment bullet with peer review time
$ws.Range("H9").Value = "- [Diagrams](.../slides/08-diagrams.qmd)`n- Peer review time – User guide"

# Match the saved view state: scrolled down/right one cell, with I10 selected
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("I10").Select()
